$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 12791.2
$ws.Range("I21").Value = 5000
$ws.Range("J21").Value = 17985.334
$ws.Range("K21").Value = 5000
$ws.Range("L21").Value = 17985.334
$ws.Range("M21").Value = -4532
$ws.Range("N21").Value = -18921.334
$ws.Range("H23").Value = 12791.2
$ws.Range("I23").Value = 5000
$ws.Range("J23").Value = 17985.334
$ws.Range("K23").Value = 5000
$ws.Range("L23").Value = 17985.334
$ws.Range("M23").Value = -4766
$ws.Range("N23").Value = -18453.334
$ws.Range("I62").Value = 27781528
$ws.Range("J62").Value = 8499.5
$ws.Range("K62").Value = 27781528
$ws.Range("L62").Value = 8499.5
$ws.Range("M62").Value = -27780904
$ws.Range("N62").Value = -9747.5
$ws.Range("I65").Value = 27781528
$ws.Range("J65").Value = 8499.5
$ws.Range("K65").Value = 138907640
$ws.Range("L65").Value = 42497.5
$ws.Range("M65").Value = -138904520
$ws.Range("N65").Value = -48737.5
$ws.Range("H80").Value = 3532
$ws.Range("I80").Value = 4843.6
$ws.Range("J80").Value = 908.8
$ws.Range("K80").Value = 14530.8
$ws.Range("L80").Value = 2726.4
$ws.Range("M80").Value = -13532.8
$ws.Range("N80").Value = -4722.4
$ws.Range("H83").Value = 3532
$ws.Range("I83").Value = 4843.6
$ws.Range("J83").Value = 908.8
$ws.Range("K83").Value = 43592.4
$ws.Range("L83").Value = 8179.2
$ws.Range("M83").Value = -38600.4
$ws.Range("N83").Value = -18163.2
$ws.Range("H96").Value = 2622.2856
$ws.Range("I96").Value = 5342.3335
$ws.Range("J96").Value = 582.25
$ws.Range("K96").Value = 16027.0005
$ws.Range("L96").Value = 1746.75
$ws.Range("M96").Value = -14654.0005
$ws.Range("N96").Value = -4492.75
$ws.Range("H108").Value = 29999
$ws.Range("J108").Value = 29999
$ws.Range("L108").Value = 29999
$ws.Range("N108").Value = -37679
$ws.Range("H116").Value = 3356
$ws.Range("I116").Value = 1973.3334
$ws.Range("J116").Value = 3874.5
$ws.Range("K116").Value = 1973.3334
$ws.Range("L116").Value = 3874.5
$ws.Range("M116").Value = 1468.6666
$ws.Range("N116").Value = -10758.5
$ws.Range("H128").Value = 36666.332
$ws.Range("J128").Value = 36666.332
$ws.Range("L128").Value = 36666.332
$ws.Range("N128").Value = -46626.332
$ws.Range("H132").Value = 7759453
$ws.Range("I132").Value = 11116693
$ws.Range("K132").Value = 33350079
$ws.Range("M132").Value = -33347549
$ws.Range("H136").Value = 32720
$ws.Range("J136").Value = 32720
$ws.Range("L136").Value = 32720
$ws.Range("N136").Value = -42920
$ws.Range("H137").Value = 1685.3715
$ws.Range("I137").Value = 1278.421
$ws.Range("J137").Value = 2168.625
$ws.Range("K137").Value = 3835.263
$ws.Range("L137").Value = 6505.875
$ws.Range("M137").Value = -1285.263
$ws.Range("N137").Value = -11605.875
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 52634030
$ws.Range("I61").Value = 58825916
$ws.Range("K61").Value = 58825916
$ws.Range("M61").Value = -58825704
$ws.Range("H74").Value = 2422.4443
$ws.Range("I74").Value = 1975.25
$ws.Range("J74").Value = 6000
$ws.Range("K74").Value = 1975.25
$ws.Range("L74").Value = 6000
$ws.Range("M74").Value = -1101.25
$ws.Range("N74").Value = -7748
$ws.Range("H77").Value = 2422.4443
$ws.Range("I77").Value = 1975.25
$ws.Range("J77").Value = 6000
$ws.Range("K77").Value = 9876.25
$ws.Range("L77").Value = 30000
$ws.Range("M77").Value = -5508.25
$ws.Range("N77").Value = -38736
$ws.Range("H97").Value = 351.90323
$ws.Range("I97").Value = 351.90323
$ws.Range("K97").Value = 351.90323
$ws.Range("M97").Value = 144.09677
$ws.Range("H136").Value = 52634030
$ws.Range("I136").Value = 58825916
$ws.Range("K136").Value = 176477748
$ws.Range("M136").Value = -176475198
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2394.3333
$ws.Range("I107").Value = 1786.6666
$ws.Range("K107").Value = 1786.6666
$ws.Range("M107").Value = 133.3334
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1311.6571
$ws.Range("I31").Value = 1162.2188
$ws.Range("J31").Value = 2905.6667
$ws.Range("K31").Value = 1162.2188
$ws.Range("L31").Value = 2905.6667
$ws.Range("M31").Value = -867.2188000000001
$ws.Range("N31").Value = -3495.6667
$ws.Range("H34").Value = 1311.6571
$ws.Range("I34").Value = 1162.2188
$ws.Range("J34").Value = 2905.6667
$ws.Range("K34").Value = 1162.2188
$ws.Range("L34").Value = 2905.6667
$ws.Range("M34").Value = -960.2188000000001
$ws.Range("N34").Value = -3309.6667
$ws.Range("H58").Value = 1462.4333
$ws.Range("I58").Value = 1039.7222
$ws.Range("J58").Value = 2096.5
$ws.Range("K58").Value = 1039.7222
$ws.Range("L58").Value = 2096.5
$ws.Range("M58").Value = -836.7221999999999
$ws.Range("N58").Value = -2502.5
$ws.Range("H99").Value = 1507.8334
$ws.Range("I99").Value = 1548.5714
$ws.Range("K99").Value = 1548.5714
$ws.Range("M99").Value = -50.57140000000004
$ws.Range("H126").Value = 1507.8334
$ws.Range("I126").Value = 1548.5714
$ws.Range("K126").Value = 4645.7142
$ws.Range("M126").Value = -2175.7142
$ws.Range("H134").Value = 166673340
$ws.Range("I134").Value = 10005.5
$ws.Range("J134").Value = 500000000
$ws.Range("K134").Value = 30016.5
$ws.Range("L134").Value = 1500000000
$ws.Range("M134").Value = -27481.5
$ws.Range("N134").Value = -1500005070
$ws.Range("H136").Value = 1462.4333
$ws.Range("I136").Value = 1039.7222
$ws.Range("J136").Value = 2096.5
$ws.Range("K136").Value = 3119.1666
$ws.Range("L136").Value = 6289.5
$ws.Range("M136").Value = -569.1665999999996
$ws.Range("N136").Value = -11389.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 100.44444
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 112.875
$ws.Range("K2").Value = 6
$ws.Range("L2").Value = 677.25
$ws.Range("M2").Value = 107
$ws.Range("N2").Value = -903.25
$ws.Range("H68").Value = 1250.3334
$ws.Range("I68").Value = 1138
$ws.Range("J68").Value = 1475
$ws.Range("K68").Value = 3414
$ws.Range("L68").Value = 4425
$ws.Range("M68").Value = -2603
$ws.Range("N68").Value = -6047
$ws.Range("H71").Value = 1250.3334
$ws.Range("I71").Value = 1138
$ws.Range("J71").Value = 1475
$ws.Range("K71").Value = 10242
$ws.Range("L71").Value = 13275
$ws.Range("M71").Value = -6186
$ws.Range("N71").Value = -21387
$ws.Range("H93").Value = 5206.727
$ws.Range("I93").Value = 1000
$ws.Range("J93").Value = 5627.4
$ws.Range("K93").Value = 3000
$ws.Range("L93").Value = 16882.2
$ws.Range("M93").Value = -1128
$ws.Range("N93").Value = -20626.2
$ws.Range("H113").Value = 688.8823
$ws.Range("I113").Value = 503
$ws.Range("J113").Value = 700.5
$ws.Range("K113").Value = 1509
$ws.Range("L113").Value = 2101.5
$ws.Range("M113").Value = 661
$ws.Range("N113").Value = -6441.5
$ws.Range("H131").Value = 20411452
$ws.Range("J131").Value = 3556.2666
$ws.Range("L131").Value = 10668.7998
$ws.Range("N131").Value = -20748.7998
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 18753204
$ws.Range("I70").Value = 17860344
$ws.Range("J70").Value = 20003208
$ws.Range("K70").Value = 17860344
$ws.Range("L70").Value = 20003208
$ws.Range("M70").Value = -17860074
$ws.Range("N70").Value = -20003748
$ws.Range("H73").Value = 18753204
$ws.Range("I73").Value = 17860344
$ws.Range("J73").Value = 20003208
$ws.Range("K73").Value = 17860344
$ws.Range("L73").Value = 20003208
$ws.Range("M73").Value = -17859408
$ws.Range("N73").Value = -20005080
$ws.Range("H100").Value = 30666.666
$ws.Range("J100").Value = 30666.666
$ws.Range("L100").Value = 30666.666
$ws.Range("N100").Value = -32830.666
$ws.Range("H141").Value = 32700
$ws.Range("J141").Value = 32700
$ws.Range("L141").Value = 32700
$ws.Range("N141").Value = -43060
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 869.1539
$ws.Range("I55").Value = 770
$ws.Range("J55").Value = 1199.6666
$ws.Range("K55").Value = 770
$ws.Range("L55").Value = 1199.6666
$ws.Range("M55").Value = -597
$ws.Range("N55").Value = -1545.6666
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3108.3333
$ws.Range("I96").Value = 2024.5
$ws.Range("J96").Value = 5276
$ws.Range("K96").Value = 2024.5
$ws.Range("L96").Value = 6100
$ws.Range("M96").Value = -651.5
$ws.Range("N96").Value = -8022
$ws.Range("H136").Value = 1176.0435
$ws.Range("I136").Value = 1184.2667
$ws.Range("J136").Value = 1160.625
$ws.Range("K136").Value = 3552.800099999999
$ws.Range("L136").Value = 3481.875
$ws.Range("M136").Value = -1002.800099999999
$ws.Range("N136").Value = -8581.875
